$d = $word.ActiveDocument

# The paragraph currently reads "My third file " (runs: "My ", "third", " ",
# "file", " "). We need it to read "My third file daffafa", split as the
# runs "My " (unchanged), "third file " (was "third"), and "daffafa" (was
# the final " ").
#
# Replace the trailing " file " (last three runs: " ", "file", " ") first so
# that the literal text "file" found by the second Find below can only match
# the still-untouched "third" run - doing the "third" replacement first
# would leave the substring "file" inside the freshly written "third file "
# run and the second Find could latch onto that instead.
$findTail = $d.Content.Find
$findTail.ClearFormatting()
$findTail.Replacement.ClearFormatting()
$findTail.Execute(" file ", $true, $false, $false, $false, $false, $true, 1, $false, "daffafa", 2)

$findHead = $d.Content.Find
$findHead.ClearFormatting()
$findHead.Replacement.ClearFormatting()
$findHead.Execute("third", $true, $false, $false, $false, $false, $true, 1, $false, "third file ", 2)
